# Update recalculated NATMI TPM-derived values for the Dkk2-Kremen1
# LR-pairs sheet ("update scripts wuth new tpm").
#
# The sending/target cluster labels, ligand/receptor symbols and overall
# row layout are unchanged; only the expression-derived metrics (columns
# E, F, G, H, I, J for the ligand side and M, N, O, P, Q, R, S, T for the
# receptor/edge side) were recomputed against the new TPM values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3447283333333333
$ws.Range("H2").Value = 1.034185
$ws.Range("I2").Value = 0.07368549602308437
$ws.Range("J2").Value = 0.07368549602308436
$ws.Range("M2").Value = 7.011769666666666
$ws.Range("N2").Value = 21.035309
$ws.Range("O2").Value = 0.1665495980465456
$ws.Range("P2").Value = 0.1665495980465456
$ws.Range("Q2").Value = 2.417155670907222
$ws.Range("R2").Value = 21.754401038165
$ws.Range("S2").Value = 0.01227228974450504
$ws.Range("T2").Value = 0.01227228974450503
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3447283333333333
$ws.Range("H3").Value = 1.034185
$ws.Range("I3").Value = 0.07368549602308437
$ws.Range("J3").Value = 0.07368549602308436
$ws.Range("O3").Value = 0.2689007917296081
$ws.Range("P3").Value = 0.2689007917296081
$ws.Range("Q3").Value = 3.902591667972777
$ws.Range("R3").Value = 35.12332501175499
$ws.Range("S3").Value = 0.01981408821959628
$ws.Range("T3").Value = 0.01981408821959628
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3447283333333333
$ws.Range("H4").Value = 1.034185
$ws.Range("I4").Value = 0.07368549602308437
$ws.Range("J4").Value = 0.07368549602308436
$ws.Range("M4").Value = 2.485142
$ws.Range("N4").Value = 7.455426000000001
$ws.Range("O4").Value = 0.05902923525229724
$ws.Range("P4").Value = 0.05902923525229725
$ws.Range("Q4").Value = 0.8566988597566667
$ws.Range("R4").Value = 7.71028973781
$ws.Range("S4").Value = 0.00434959847942886
$ws.Range("T4").Value = 0.00434959847942886
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.3447283333333333
$ws.Range("H5").Value = 1.034185
$ws.Range("I5").Value = 0.07368549602308437
$ws.Range("J5").Value = 0.07368549602308436
$ws.Range("M5").Value = 19.36045366666667
$ws.Range("N5").Value = 58.081361
$ws.Range("O5").Value = 0.4598661863510687
$ws.Range("P5").Value = 0.4598661863510686
$ws.Range("Q5").Value = 6.674096925087222
$ws.Range("R5").Value = 60.066872325785
$ws.Range("S5").Value = 0.03388546804552265
$ws.Range("T5").Value = 0.03388546804552264
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.3447283333333333
$ws.Range("H6").Value = 1.034185
$ws.Range("I6").Value = 0.07368549602308437
$ws.Range("J6").Value = 0.07368549602308436
$ws.Range("M6").Value = 1.92205
$ws.Range("N6").Value = 5.76615
$ws.Range("O6").Value = 0.0456541886204804
$ws.Range("P6").Value = 0.0456541886204804
$ws.Range("Q6").Value = 0.6625850930833332
$ws.Range("R6").Value = 5.963265837749999
$ws.Range("S6").Value = 0.003364051534031552
$ws.Range("T6").Value = 0.003364051534031552
$ws.Range("I7").Value = 0.905070728628456
$ws.Range("J7").Value = 0.9050707286284558
$ws.Range("M7").Value = 7.011769666666666
$ws.Range("N7").Value = 21.035309
$ws.Range("O7").Value = 0.1665495980465456
$ws.Range("P7").Value = 0.1665495980465456
$ws.Range("Q7").Value = 29.68965349152344
$ws.Range("R7").Value = 267.206881423711
$ws.Range("S7").Value = 0.1507391660567635
$ws.Range("T7").Value = 0.1507391660567635
$ws.Range("I8").Value = 0.905070728628456
$ws.Range("J8").Value = 0.9050707286284558
$ws.Range("O8").Value = 0.2689007917296081
$ws.Range("P8").Value = 0.2689007917296081
$ws.Range("S8").Value = 0.2433742354994851
$ws.Range("T8").Value = 0.2433742354994851
$ws.Range("I9").Value = 0.905070728628456
$ws.Range("J9").Value = 0.9050707286284558
$ws.Range("M9").Value = 2.485142
$ws.Range("N9").Value = 7.455426000000001
$ws.Range("O9").Value = 0.05902923525229724
$ws.Range("P9").Value = 0.05902923525229725
$ws.Range("Q9").Value = 10.52273653653933
$ws.Range("R9").Value = 94.70462882885401
$ws.Range("S9").Value = 0.0534256329601772
$ws.Range("T9").Value = 0.0534256329601772
$ws.Range("I10").Value = 0.905070728628456
$ws.Range("J10").Value = 0.9050707286284558
$ws.Range("M10").Value = 19.36045366666667
$ws.Range("N10").Value = 58.081361
$ws.Range("O10").Value = 0.4598661863510687
$ws.Range("P10").Value = 0.4598661863510686
$ws.Range("Q10").Value = 81.97718808913545
$ws.Range("R10").Value = 737.794692802219
$ws.Range("S10").Value = 0.416211424352351
$ws.Range("T10").Value = 0.4162114243523509
$ws.Range("I11").Value = 0.905070728628456
$ws.Range("J11").Value = 0.9050707286284558
$ws.Range("M11").Value = 1.92205
$ws.Range("N11").Value = 5.76615
$ws.Range("O11").Value = 0.0456541886204804
$ws.Range("P11").Value = 0.0456541886204804
$ws.Range("Q11").Value = 8.138458792316666
$ws.Range("R11").Value = 73.24612913085
$ws.Range("S11").Value = 0.04132026975967916
$ws.Range("T11").Value = 0.04132026975967915
$ws.Range("G12").Value = 0.027522
$ws.Range("H12").Value = 0.082566
$ws.Range("I12").Value = 0.005882812712079546
$ws.Range("J12").Value = 0.005882812712079545
$ws.Range("M12").Value = 7.011769666666666
$ws.Range("N12").Value = 21.035309
$ws.Range("O12").Value = 0.1665495980465456
$ws.Range("P12").Value = 0.1665495980465456
$ws.Range("Q12").Value = 0.192977924766
$ws.Range("R12").Value = 1.736801322894
$ws.Range("S12").Value = 0.0009797800925799571
$ws.Range("T12").Value = 0.000979780092579957
$ws.Range("G13").Value = 0.027522
$ws.Range("H13").Value = 0.082566
$ws.Range("I13").Value = 0.005882812712079546
$ws.Range("J13").Value = 0.005882812712079545
$ws.Range("O13").Value = 0.2689007917296081
$ws.Range("P13").Value = 0.2689007917296081
$ws.Range("Q13").Value = 0.311570351202
$ws.Range("R13").Value = 2.804133160818
$ws.Range("S13").Value = 0.001581892995875193
$ws.Range("T13").Value = 0.001581892995875193
$ws.Range("G14").Value = 0.027522
$ws.Range("H14").Value = 0.082566
$ws.Range("I14").Value = 0.005882812712079546
$ws.Range("J14").Value = 0.005882812712079545
$ws.Range("M14").Value = 2.485142
$ws.Range("N14").Value = 7.455426000000001
$ws.Range("O14").Value = 0.05902923525229724
$ws.Range("P14").Value = 0.05902923525229725
$ws.Range("Q14").Value = 0.068396078124
$ws.Range("R14").Value = 0.6155647031160001
$ws.Range("S14").Value = 0.0003472579355265483
$ws.Range("T14").Value = 0.0003472579355265483
$ws.Range("G15").Value = 0.027522
$ws.Range("H15").Value = 0.082566
$ws.Range("I15").Value = 0.005882812712079546
$ws.Range("J15").Value = 0.005882812712079545
$ws.Range("M15").Value = 19.36045366666667
$ws.Range("N15").Value = 58.081361
$ws.Range("O15").Value = 0.4598661863510687
$ws.Range("P15").Value = 0.4598661863510686
$ws.Range("Q15").Value = 0.5328384058140001
$ws.Range("R15").Value = 4.795545652326
$ws.Range("S15").Value = 0.002705306646921608
$ws.Range("T15").Value = 0.002705306646921608
$ws.Range("G16").Value = 0.027522
$ws.Range("H16").Value = 0.082566
$ws.Range("I16").Value = 0.005882812712079546
$ws.Range("J16").Value = 0.005882812712079545
$ws.Range("M16").Value = 1.92205
$ws.Range("N16").Value = 5.76615
$ws.Range("O16").Value = 0.0456541886204804
$ws.Range("P16").Value = 0.0456541886204804
$ws.Range("Q16").Value = 0.05289866009999999
$ws.Range("R16").Value = 0.4760879409
$ws.Range("S16").Value = 0.0002685750411762395
$ws.Range("T16").Value = 0.0002685750411762394
$ws.Range("G17").Value = 0.07186433333333334
$ws.Range("H17").Value = 0.215593
$ws.Range("I17").Value = 0.01536096263638017
$ws.Range("J17").Value = 0.01536096263638017
$ws.Range("M17").Value = 7.011769666666666
$ws.Range("N17").Value = 21.035309
$ws.Range("O17").Value = 0.1665495980465456
$ws.Range("P17").Value = 0.1665495980465456
$ws.Range("Q17").Value = 0.5038961525818889
$ws.Range("R17").Value = 4.535065373237
$ws.Range("S17").Value = 0.002558362152697123
$ws.Range("T17").Value = 0.002558362152697123
$ws.Range("G18").Value = 0.07186433333333334
$ws.Range("H18").Value = 0.215593
$ws.Range("I18").Value = 0.01536096263638017
$ws.Range("J18").Value = 0.01536096263638017
$ws.Range("O18").Value = 0.2689007917296081
$ws.Range("P18").Value = 0.2689007917296081
$ws.Range("Q18").Value = 0.8135599002821111
$ws.Range("R18").Value = 7.322039102539
$ws.Range("S18").Value = 0.004130575014651557
$ws.Range("T18").Value = 0.004130575014651557
$ws.Range("G19").Value = 0.07186433333333334
$ws.Range("H19").Value = 0.215593
$ws.Range("I19").Value = 0.01536096263638017
$ws.Range("J19").Value = 0.01536096263638017
$ws.Range("M19").Value = 2.485142
$ws.Range("N19").Value = 7.455426000000001
$ws.Range("O19").Value = 0.05902923525229724
$ws.Range("P19").Value = 0.05902923525229725
$ws.Range("Q19").Value = 0.1785930730686667
$ws.Range("R19").Value = 1.607337657618
$ws.Range("S19").Value = 0.0009067458771646334
$ws.Range("T19").Value = 0.0009067458771646335
$ws.Range("G20").Value = 0.07186433333333334
$ws.Range("H20").Value = 0.215593
$ws.Range("I20").Value = 0.01536096263638017
$ws.Range("J20").Value = 0.01536096263638017
$ws.Range("M20").Value = 19.36045366666667
$ws.Range("N20").Value = 58.081361
$ws.Range("O20").Value = 0.4598661863510687
$ws.Range("P20").Value = 0.4598661863510686
$ws.Range("Q20").Value = 1.391326095785889
$ws.Range("R20").Value = 12.521934862073
$ws.Range("S20").Value = 0.007063987306273408
$ws.Range("T20").Value = 0.007063987306273408
$ws.Range("G21").Value = 0.07186433333333334
$ws.Range("H21").Value = 0.215593
$ws.Range("I21").Value = 0.01536096263638017
$ws.Range("J21").Value = 0.01536096263638017
$ws.Range("M21").Value = 1.92205
$ws.Range("N21").Value = 5.76615
$ws.Range("O21").Value = 0.0456541886204804
$ws.Range("P21").Value = 0.0456541886204804
$ws.Range("Q21").Value = 0.1381268418833333
$ws.Range("R21").Value = 1.24314157695
$ws.Range("S21").Value = 0.0007012922855934524
$ws.Range("T21").Value = 0.0007012922855934524
